$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.401.36"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.802.34"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'227.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'0.581"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'35.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.37%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'0.0691"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "2.063.19"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "1.807.68"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "34.399.13"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "'4.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").Value = "'68.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0795"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'245.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").Value = "'11.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "'4.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").Value = "'170.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("E26").Value = "  +3.91%  "
$ws.Range("D27").Value = "'16.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "'0.119"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'3.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "'3.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "1.396.59"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("D36").Value = "'0.677"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'2.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").Value = "'82.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("D42").Value = "'0.946"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "'13.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("D47").Value = "'5.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "1.963.32"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "'104.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "0.0₆0129"
$ws.Range("E51").Value = "  +0.07%  "
